$d = $word.ActiveDocument

# --- Add the new paragraph at the end of the document body ---------------
$body = '<w:p><w:r><w:t>OK, here is a new addition to our experiment</w:t></w:r><w:r><w:t>!! Ha-ha-ha.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Thomaz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$c = $d.Content
$c.Collapse(0)
$c.InsertXML($xml)

# --- sectPr tweaks: header/footer distance + column spacing --------------
$ps = $d.PageSetup
$ps.HeaderDistance = 36
$ps.FooterDistance = 36
$ps.TextColumns.Spacing = 36
